$d = $word.ActiveDocument

# Locate the last stanza's paragraph ("И падаю в забвение экстаза") robustly,
# without assuming a fixed paragraph index.
$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains("падаю в забвение")) {
        $target = $p
    }
}

$startPos = $target.Range.Start
$endPos = $d.Content.End
$replaceRange = $d.Range($startPos, $endPos)

# Rebuild that paragraph (giving the paragraph mark and the final run
# "экстаза" the same white/background1 color already used by the other
# runs) and append the two new white-text stanzas:
#   "Не против жестокого проказа,"
#   "Ведь мы живем лишь пару " + "лет" (last word left uncolored)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:color w:val="FFFFFF" w:themeColor="background1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FFFFFF" w:themeColor="background1"/></w:rPr><w:t xml:space="preserve">И </w:t></w:r><w:r><w:rPr><w:color w:val="FFFFFF" w:themeColor="background1"/></w:rPr><w:t xml:space="preserve">падаю в забвение </w:t></w:r><w:r><w:rPr><w:color w:val="FFFFFF" w:themeColor="background1"/></w:rPr><w:t>экстаза</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="FFFFFF" w:themeColor="background1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FFFFFF" w:themeColor="background1"/></w:rPr><w:t>Не против жестокого проказа,</w:t></w:r></w:p><w:p><w:r><w:rPr><w:color w:val="FFFFFF" w:themeColor="background1"/></w:rPr><w:t xml:space="preserve">Ведь мы живем лишь пару </w:t></w:r><w:r><w:t>лет</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$replaceRange.InsertXML($xml)
